$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8807.706
$ws.Range("I86").Value = 3003.6
$ws.Range("J86").Value = 17099.285
$ws.Range("K86").Value = 3003.6
$ws.Range("L86").Value = 17099.285
$ws.Range("M86").Value = -1880.6
$ws.Range("N86").Value = -19345.285
$ws.Range("H89").Value = 8807.706
$ws.Range("I89").Value = 3003.6
$ws.Range("J89").Value = 17099.285
$ws.Range("K89").Value = 15018
$ws.Range("L89").Value = 85496.425
$ws.Range("M89").Value = -9402
$ws.Range("N89").Value = -96728.425
$ws.Range("H99").Value = 167.25
$ws.Range("I99").Value = 167.25
$ws.Range("K99").Value = 501.75
$ws.Range("M99").Value = 996.25
$ws.Range("H113").Value = 37040828
$ws.Range("I113").Value = 71431770
$ws.Range("J113").Value = 4427.5386
$ws.Range("K113").Value = 71431770
$ws.Range("L113").Value = 4427.5386
$ws.Range("M113").Value = -71428516
$ws.Range("N113").Value = -10935.5386
$ws.Range("H129").Value = 589661.5600000001
$ws.Range("J129").Value = 771004.4
$ws.Range("L129").Value = 2313013.2
$ws.Range("N129").Value = -2323013.2
$ws.Range("H138").Value = 2186.8096
$ws.Range("I138").Value = 1698.7368
$ws.Range("J138").Value = 2590
$ws.Range("K138").Value = 5096.2104
$ws.Range("L138").Value = 7770
$ws.Range("M138").Value = 43.78960000000006
$ws.Range("N138").Value = -18050
$ws.Range("H141").Value = 1243.1316
$ws.Range("I141").Value = 954.1177
$ws.Range("K141").Value = 2862.3531
$ws.Range("M141").Value = 2317.6469

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1612.3158
$ws.Range("I2").Value = 1470
$ws.Range("K2").Value = 1470
$ws.Range("M2").Value = -1357
$ws.Range("H32").Value = 5268.1055
$ws.Range("I32").Value = 5183.137
$ws.Range("K32").Value = 5183.137
$ws.Range("M32").Value = -4896.137
$ws.Range("H45").Value = 2750
$ws.Range("I45").Value = 2573.9048
$ws.Range("J45").Value = 3212.25
$ws.Range("K45").Value = 2573.9048
$ws.Range("L45").Value = 3212.25
$ws.Range("M45").Value = -2196.9048
$ws.Range("N45").Value = -3966.25
$ws.Range("H102").Value = 1269.75
$ws.Range("I102").Value = 789.5
$ws.Range("K102").Value = 789.5
$ws.Range("M102").Value = 832.5
$ws.Range("H110").Value = 749.5
$ws.Range("I110").Value = 749.5
$ws.Range("K110").Value = 749.5
$ws.Range("M110").Value = 1295.5
$ws.Range("H116").Value = 1612.3158
$ws.Range("I116").Value = 1470
$ws.Range("K116").Value = 1470
$ws.Range("M116").Value = 824
$ws.Range("H132").Value = 28108.018
$ws.Range("I132").Value = 1538.4667
$ws.Range("J132").Value = 127743.836
$ws.Range("K132").Value = 4615.4001
$ws.Range("L132").Value = 383231.508
$ws.Range("M132").Value = -2085.4001
$ws.Range("N132").Value = -388291.508

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1612.3158
$ws.Range("I3").Value = 1470
$ws.Range("K3").Value = 1470
$ws.Range("M3").Value = -1356
$ws.Range("H105").Value = 4213.8
$ws.Range("I105").Value = 4806.3335
$ws.Range("K105").Value = 4806.3335
$ws.Range("M105").Value = -3059.3335
$ws.Range("H107").Value = 669.55
$ws.Range("I107").Value = 485.55554
$ws.Range("K107").Value = 485.55554
$ws.Range("M107").Value = 1434.44446

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 39749.5
$ws.Range("J52").Value = 39749.5
$ws.Range("L52").Value = 39749.5
$ws.Range("N52").Value = -40337.5
$ws.Range("H132").Value = 13937.714
$ws.Range("I132").Value = 14820.486
$ws.Range("K132").Value = 44461.458
$ws.Range("M132").Value = -41931.458
$ws.Range("H134").Value = 642.85364
$ws.Range("I134").Value = 528.7353000000001
$ws.Range("J134").Value = 1197.1428
$ws.Range("K134").Value = 1586.2059
$ws.Range("L134").Value = 3591.4284
$ws.Range("M134").Value = 948.7940999999998
$ws.Range("N134").Value = -8661.428400000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1911.7778
$ws.Range("I5").Value = 1350.875
$ws.Range("J5").Value = 6399
$ws.Range("K5").Value = 4052.625
$ws.Range("L5").Value = 19197
$ws.Range("M5").Value = -3940.625
$ws.Range("N5").Value = -19421
$ws.Range("H25").Value = 724.75
$ws.Range("I25").Value = 724.75
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 2174.25
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -2005.25
$ws.Range("N25").ClearContents()
$ws.Range("H30").Value = 724.75
$ws.Range("I30").Value = 724.75
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2174.25
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -2072.25
$ws.Range("N30").ClearContents()
$ws.Range("H122").Value = 745.875
$ws.Range("I122").Value = 353
$ws.Range("J122").Value = 942.3125
$ws.Range("K122").Value = 3177
$ws.Range("L122").Value = 8480.8125
$ws.Range("M122").Value = -727
$ws.Range("N122").Value = -13380.8125
$ws.Range("H131").Value = 164748.73
$ws.Range("J131").Value = 176276.19
$ws.Range("L131").Value = 528828.5700000001
$ws.Range("N131").Value = -538908.5700000001
$ws.Range("H132").Value = 963.3570999999999
$ws.Range("I132").Value = 479.6
$ws.Range("K132").Value = 4316.400000000001
$ws.Range("M132").Value = -1786.400000000001
$ws.Range("H135").Value = 1911.7778
$ws.Range("I135").Value = 1350.875
$ws.Range("J135").Value = 6399
$ws.Range("K135").Value = 12157.875
$ws.Range("L135").Value = 57591
$ws.Range("M135").Value = -9622.875
$ws.Range("N135").Value = -62661

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 70.64286
$ws.Range("I2").Value = 62.375
$ws.Range("K2").Value = 62.375
$ws.Range("M2").Value = 50.625
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H43").Value = 188181.81
$ws.Range("I43").Value = 2000000
$ws.Range("K43").Value = 2000000
$ws.Range("M43").Value = -1999849
$ws.Range("H46").Value = 20596
$ws.Range("J46").Value = 20596
$ws.Range("L46").Value = 20596
$ws.Range("N46").Value = -20908
$ws.Range("H57").Value = 25540
$ws.Range("J57").Value = 24920
$ws.Range("L57").Value = 24920
$ws.Range("N57").Value = -26560
$ws.Range("H80").Value = 4310.6
$ws.Range("J80").Value = 4529.4287
$ws.Range("L80").Value = 4529.4287
$ws.Range("N80").Value = -6525.4287
$ws.Range("H83").Value = 4310.6
$ws.Range("J83").Value = 4529.4287
$ws.Range("L83").Value = 22647.1435
$ws.Range("N83").Value = -32631.1435
$ws.Range("H102").Value = 38464884
$ws.Range("I102").Value = 38464884
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 38464884
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -38463262
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 78432730
$ws.Range("I122").Value = 33334192
$ws.Range("K122").Value = 100002576
$ws.Range("M122").Value = -100000126
$ws.Range("H132").Value = 16911.555
$ws.Range("I132").Value = 3145.889
$ws.Range("K132").Value = 9437.667000000001
$ws.Range("M132").Value = -6907.667000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3078.923
$ws.Range("I40").Value = 2402.625
$ws.Range("K40").Value = 2402.625
$ws.Range("M40").Value = -2266.625
$ws.Range("H61").Value = 5442.25
$ws.Range("I61").Value = 2388.375
$ws.Range("K61").Value = 2388.375
$ws.Range("M61").Value = -2186.375
$ws.Range("H113").Value = 5442.25
$ws.Range("I113").Value = 2388.375
$ws.Range("K113").Value = 2388.375
$ws.Range("M113").Value = -218.375
$ws.Range("H132").Value = 1881.1538
$ws.Range("I132").Value = 1539.6522
$ws.Range("J132").Value = 4499.3335
$ws.Range("K132").Value = 4618.9566
$ws.Range("L132").Value = 13498.0005
$ws.Range("M132").Value = -2088.9566
$ws.Range("N132").Value = -18558.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 1999
$ws.Range("I19").Value = 1999
$ws.Range("K19").Value = 1999
$ws.Range("M19").Value = -1825
$ws.Range("H113").Value = 2254772.5
$ws.Range("I113").Value = 3454.875
$ws.Range("K113").Value = 10364.625
$ws.Range("M113").Value = -8194.625
$ws.Range("H122").Value = 1175.5714
$ws.Range("I122").Value = 1079
$ws.Range("J122").Value = 1263.3636
$ws.Range("K122").Value = 3237
$ws.Range("L122").Value = 3790.0908
$ws.Range("M122").Value = -787
$ws.Range("N122").Value = -8690.0908
$ws.Range("H132").Value = 922.80646
$ws.Range("I132").Value = 641.6896400000001
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 1925.06892
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = 604.9310799999998
$ws.Range("N132").Value = -20057
